# Auto-generated edit script: update Ifrit_Profits market-data snapshot values
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("I31").Value = 800
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2400
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2170
$ws.Range("N31").ClearContents()

# Row 32
$ws.Range("H32").Value = 6722.5
$ws.Range("I32").Value = 7445
$ws.Range("J32").Value = 6000
$ws.Range("K32").Value = 7445
$ws.Range("L32").Value = 6000
$ws.Range("M32").Value = -7119
$ws.Range("N32").Value = -6652

# Row 33
$ws.Range("H33").Value = 192.72223
$ws.Range("I33").Value = 192.72223
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 192.72223
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 36.27777
$ws.Range("N33").ClearContents()

# Row 64
$ws.Range("H64").Value = 4236.4116
$ws.Range("I64").Value = 4222
$ws.Range("J64").Value = 4257
$ws.Range("K64").Value = 4222
$ws.Range("L64").Value = 4257
$ws.Range("M64").Value = -3974
$ws.Range("N64").Value = -4753

# Row 67
$ws.Range("H67").Value = 4236.4116
$ws.Range("I67").Value = 4222
$ws.Range("J67").Value = 4257
$ws.Range("K67").Value = 4222
$ws.Range("L67").Value = 4257
$ws.Range("M67").Value = -3364
$ws.Range("N67").Value = -5973

# Row 70
$ws.Range("H70").Value = 1856.5
$ws.Range("I70").Value = 1738
$ws.Range("J70").Value = 1975
$ws.Range("K70").Value = 5214
$ws.Range("L70").Value = 5925
$ws.Range("M70").Value = -4944
$ws.Range("N70").Value = -6465

# Row 73
$ws.Range("H73").Value = 1856.5
$ws.Range("I73").Value = 1738
$ws.Range("J73").Value = 1975
$ws.Range("K73").Value = 5214
$ws.Range("L73").Value = 5925
$ws.Range("M73").Value = -4278
$ws.Range("N73").Value = -7797

# Row 76
$ws.Range("H76").Value = 3079.3333
$ws.Range("I76").Value = 3050
$ws.Range("J76").Value = 3490
$ws.Range("K76").Value = 3050
$ws.Range("L76").Value = 3490
$ws.Range("M76").Value = -2735
$ws.Range("N76").Value = -4120

# Row 79
$ws.Range("H79").Value = 3079.3333
$ws.Range("I79").Value = 3050
$ws.Range("J79").Value = 3490
$ws.Range("K79").Value = 3050
$ws.Range("L79").Value = 3490
$ws.Range("M79").Value = -1958
$ws.Range("N79").Value = -5674

# Row 98
$ws.Range("H98").Value = 2266.475
$ws.Range("I98").Value = 2404.25
$ws.Range("J98").Value = 1026.5
$ws.Range("K98").Value = 2404.25
$ws.Range("L98").Value = 1026.5
$ws.Range("M98").Value = -906.25
$ws.Range("N98").Value = -4022.5

# Row 107
$ws.Range("H107").Value = 842.65515
$ws.Range("I107").Value = 868.3158
$ws.Range("J107").Value = 793.9
$ws.Range("K107").Value = 868.3158
$ws.Range("L107").Value = 793.9
$ws.Range("M107").Value = 1051.6842
$ws.Range("N107").Value = -4633.9

# Row 114
$ws.Range("H114").Value = 44583.855
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 44583.855
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 44583.855
$ws.Range("N114").Value = -53261.855

# Row 122
$ws.Range("H122").Value = 2266.475
$ws.Range("I122").Value = 2404.25
$ws.Range("J122").Value = 1026.5
$ws.Range("K122").Value = 7212.75
$ws.Range("L122").Value = 3079.5
$ws.Range("M122").Value = -4762.75
$ws.Range("N122").Value = -7979.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4924.38
$ws.Range("I32").Value = 4468.4536
$ws.Range("J32").Value = 19666
$ws.Range("K32").Value = 4468.4536
$ws.Range("L32").Value = 19666
$ws.Range("M32").Value = -4181.4536
$ws.Range("N32").Value = -20240

# Row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 132
$ws.Range("H132").Value = 988666.6
$ws.Range("I132").Value = 1115229.8
$ws.Range("J132").Value = 203975.6
$ws.Range("K132").Value = 3345689.4
$ws.Range("L132").Value = 611926.8
$ws.Range("M132").Value = -3343159.4
$ws.Range("N132").Value = -616986.8

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1701
$ws.Range("I20").Value = 1953.7333
$ws.Range("J20").Value = 1279.7778
$ws.Range("K20").Value = 1953.7333
$ws.Range("L20").Value = 1279.7778
$ws.Range("M20").Value = -1706.7333
$ws.Range("N20").Value = -1773.7778

# Row 86
$ws.Range("H86").Value = 928.5454999999999
$ws.Range("I86").Value = 874.73334
$ws.Range("J86").Value = 1170.7
$ws.Range("K86").Value = 874.73334
$ws.Range("L86").Value = 1170.7
$ws.Range("M86").Value = 248.26666
$ws.Range("N86").Value = -3416.7

# Row 89
$ws.Range("H89").Value = 928.5454999999999
$ws.Range("I89").Value = 874.73334
$ws.Range("J89").Value = 1170.7
$ws.Range("K89").Value = 4373.6667
$ws.Range("L89").Value = 5853.5
$ws.Range("M89").Value = 1242.3333
$ws.Range("N89").Value = -17085.5

# Row 105
$ws.Range("H105").Value = 1540.375
$ws.Range("I105").Value = 1466.7894
$ws.Range("J105").Value = 1820
$ws.Range("K105").Value = 1466.7894
$ws.Range("L105").Value = 1820
$ws.Range("M105").Value = 280.2106000000001
$ws.Range("N105").Value = -5314

# Row 134
$ws.Range("H134").Value = 11830528
$ws.Range("I134").Value = 15470243
$ws.Range("J134").Value = 1456
$ws.Range("K134").Value = 46410729
$ws.Range("L134").Value = 4368
$ws.Range("M134").Value = -46408194
$ws.Range("N134").Value = -9438

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 2466.95
$ws.Range("I62").Value = 2276.6
$ws.Range("J62").Value = 3038
$ws.Range("K62").Value = 2276.6
$ws.Range("L62").Value = 3038
$ws.Range("M62").Value = -1652.6
$ws.Range("N62").Value = -4286

# Row 65
$ws.Range("H65").Value = 2466.95
$ws.Range("I65").Value = 2276.6
$ws.Range("J65").Value = 3038
$ws.Range("K65").Value = 11383
$ws.Range("L65").Value = 15190
$ws.Range("M65").Value = -8263
$ws.Range("N65").Value = -21430

# Row 110
$ws.Range("H110").Value = 48892
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 48892
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 48892
$ws.Range("N110").Value = -57072

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 11905483
$ws.Range("I5").Value = 593.13336
$ws.Range("J5").Value = 18519310
$ws.Range("K5").Value = 1779.40008
$ws.Range("L5").Value = 55557930
$ws.Range("M5").Value = -1667.40008
$ws.Range("N5").Value = -55558154

# Row 122
$ws.Range("H122").Value = 12319846
$ws.Range("I122").Value = 18519080
$ws.Range("J122").Value = 2175643.8
$ws.Range("K122").Value = 166671720
$ws.Range("L122").Value = 19580794.2
$ws.Range("M122").Value = -166669270
$ws.Range("N122").Value = -19585694.2

# Row 131
$ws.Range("H131").Value = 7766.3228
$ws.Range("I131").Value = 10421.667
$ws.Range("J131").Value = 7129.04
$ws.Range("K131").Value = 31265.001
$ws.Range("L131").Value = 21387.12
$ws.Range("M131").Value = -26225.001
$ws.Range("N131").Value = -31467.12

# Row 135
$ws.Range("H135").Value = 11905483
$ws.Range("I135").Value = 593.13336
$ws.Range("J135").Value = 18519310
$ws.Range("K135").Value = 5338.20024
$ws.Range("L135").Value = 166673790
$ws.Range("M135").Value = -2803.20024
$ws.Range("N135").Value = -166678860

$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Range("H32").Value = 25260
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 25260
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 25260
$ws.Range("N32").Value = -25852

# Row 97
$ws.Range("H97").Value = 2372.2222
$ws.Range("I97").Value = 2372.2222
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2372.2222
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1876.2222
$ws.Range("N97").ClearContents()

# Row 113
$ws.Range("H113").Value = 1517.2858
$ws.Range("I113").Value = 1273.5
$ws.Range("J113").Value = 2980
$ws.Range("K113").Value = 1273.5
$ws.Range("L113").Value = 2980
$ws.Range("M113").Value = 896.5

# Row 122
$ws.Range("H122").Value = 5076.4287
$ws.Range("I122").Value = 5714.1665
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 17142.4995
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -14692.4995
$ws.Range("N122").Value = -8650

$ws = $wb.Worksheets.Item("LTW")
# Row 127
$ws.Range("H127").Value = 28706.428
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 28706.428
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 28706.428
$ws.Range("N127").Value = -38626.428

# Row 139
$ws.Range("H139").Value = 75548.336
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 75548.336
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 75548.336
$ws.Range("N139").Value = -85828.336

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 14998
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 14998
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 14998
$ws.Range("N40").Value = -15296

# Row 132
$ws.Range("H132").Value = 4247.479
$ws.Range("I132").Value = 4836.3423
$ws.Range("J132").Value = 2009.8
$ws.Range("K132").Value = 14509.0269
$ws.Range("L132").Value = 6029.4
$ws.Range("M132").Value = -11979.0269
$ws.Range("N132").Value = -11089.4
